$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 885 through 1149 (inclusive) contain "Population" records whose
# year (column B) needs to move from 2022 to 2023.
$ws.Range("B885:B1149").Value = 2023
